# Applies the Feb-12-2023 GitHub Actions "Updated symbol list" refresh:
# new Price (column D) and Volume(1h) (column E) quotes for each coin row.
#
# The Price/Volume columns are stored as plain text (e.g. "309.45", "0.36%")
# rather than numbers, so every cell is momentarily switched to the "@" text
# number format before the value is written (otherwise Excel would silently
# coerce the numeric-looking strings into real numbers/percentages), and then
# restored to the workbook's default "Normal" style so no visible formatting
# change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "309.45"
Set-TextValue "E2" "0.36%"
Set-TextValue "D3" "41.18"
Set-TextValue "E3" "0.84%"
Set-TextValue "D4" "5.218"
Set-TextValue "E4" "2.21%"
Set-TextValue "D5" "0.07688"
Set-TextValue "E5" "0.89%"
Set-TextValue "D6" "1.642"
Set-TextValue "E6" "2.22%"
Set-TextValue "D7" "0.9154"
Set-TextValue "E7" "1.47%"
Set-TextValue "E8" "-1.61%"
Set-TextValue "D9" "0.1248"
Set-TextValue "E9" "10.89%"
Set-TextValue "D10" "0.1825"
Set-TextValue "E10" "1.97%"
Set-TextValue "D11" "0.09169"
Set-TextValue "E11" "0.38%"
Set-TextValue "D12" "0.04218"
Set-TextValue "E12" "0.10%"
Set-TextValue "E13" "0.03%"
Set-TextValue "D14" "0.001253"
Set-TextValue "E14" "-0.54%"
Set-TextValue "D15" "0.005746"
Set-TextValue "E15" "1.47%"
Set-TextValue "D16" "3.351"
Set-TextValue "E16" "0.09%"
Set-TextValue "D17" "4.314"
Set-TextValue "E17" "1.37%"
Set-TextValue "E18" "1.30%"
Set-TextValue "D19" "7.417"
Set-TextValue "E19" "11.83%"
Set-TextValue "D20" "0.1402"
Set-TextValue "E20" "2.67%"
Set-TextValue "D21" "0.2821"
Set-TextValue "E21" "0.68%"
Set-TextValue "D22" "0.04036"
Set-TextValue "E22" "-0.99%"
Set-TextValue "D23" "0.001266"
Set-TextValue "E23" "1.66%"
Set-TextValue "E24" "-0.15%"
Set-TextValue "D25" "0.0001302"
Set-TextValue "E25" "0.07%"
Set-TextValue "D38" "0.02569"
Set-TextValue "E38" "7.63%"
Set-TextValue "D39" "0.05351"
Set-TextValue "E39" "3.21%"
Set-TextValue "D40" "0.007847"
Set-TextValue "E40" "0.90%"
Set-TextValue "E41" "1.18%"
Set-TextValue "D42" "0.006667"
Set-TextValue "E42" "-5.57%"
Set-TextValue "D43" "0.001863"
Set-TextValue "E43" "-4.54%"
Set-TextValue "D44" "0.008064"
Set-TextValue "E44" "4.33%"
Set-TextValue "D45" "0.3071"
Set-TextValue "E45" "-0.26%"
Set-TextValue "D46" "0.00006718"
Set-TextValue "E46" "-3.65%"
Set-TextValue "E47" "0.04%"
Set-TextValue "D48" "0.2830"
Set-TextValue "E48" "511.52%"
Set-TextValue "D50" "0.00002103"
Set-TextValue "E50" "0.04%"
Set-TextValue "D51" "0.0002003"
Set-TextValue "E51" "0.04%"
